$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix status name / labels across the sheet (exact, whole-cell matches)
$xlWhole = 1          # xlWhole
$xlByRows = 1         # xlByRows
$xlNext = 1           # xlNext

$ws.Cells.Replace("bleu", "noir", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $xlWhole, $xlByRows, $false, $false, $false)
